$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 06:05"

# Pakistan (row 24) - updated case counts
$ws.Range("B24").Value = 314616
$ws.Range("C24").Value = 632
$ws.Range("D24").Value = 298968
$ws.Range("E24").Value = 9135
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 6513

# Belgica overtakes Catar in total cases, so they swap ranking order.
# Row 35 now holds Belgica (with fresh data), row 36 now holds Catar
# (with Belgica's former data, unchanged).
$ws.Range("A35").Value = "Belgica"
$ws.Range("B35").Value = 127623
$ws.Range("C35").Value = 3389
$ws.Range("D35").Value = 19645
$ws.Range("E35").Value = 97934
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 10044

$ws.Range("A36").Value = "Catar"
$ws.Range("B36").Value = 126339
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 123302
$ws.Range("E36").Value = 2821
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 216

# Honduras (row 53)
$ws.Range("B53").Value = 78788
$ws.Range("C53").Value = 519
$ws.Range("D53").Value = 29187
$ws.Range("E53").Value = 47202
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 2399

# Venezuela (row 56)
$ws.Range("B56").Value = 77646
$ws.Range("D56").Value = 68098
$ws.Range("E56").Value = 8899
$ws.Range("H56").Value = 649

# Tailandia (row 141)
$ws.Range("B141").Value = 3585
$ws.Range("C141").Value = 2
$ws.Range("D141").Value = 3388

# Burkina Faso (row 153)
$ws.Range("B153").Value = 2154
$ws.Range("D153").Value = 1397
$ws.Range("E153").Value = 698

# Butan (row 187)
$ws.Range("D187").Value = 230
$ws.Range("E187").Value = 53
